$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cx3cl1"
$ws.Range("C2").Value = "Cx3cr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 16.6160005
$ws.Range("H2").Value = 33.232001
$ws.Range("I2").Value = 0.6672963354196896
$ws.Range("J2").Value = 0.6022128312718646
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.2438225
$ws.Range("N2").Value = 0.487645
$ws.Range("O2").Value = 0.005769507646004085
$ws.Range("P2").Value = 0.003853749843732457
$ws.Range("Q2").Value = 4.051354781911249
$ws.Range("R2").Value = 16.205419127645
$ws.Range("S2").Value = 0.003849971309354406
$ws.Range("T2").Value = 0.002320777604407629

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cx3cl1"
$ws.Range("C3").Value = "Cx3cr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 16.6160005
$ws.Range("H3").Value = 33.232001
$ws.Range("I3").Value = 0.6672963354196896
$ws.Range("J3").Value = 0.6022128312718646
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1619216666666667
$ws.Range("N3").Value = 0.485765
$ws.Range("O3").Value = 0.003831509782268077
$ws.Range("P3").Value = 0.003838892622380414
$ws.Range("Q3").Value = 2.690490494294166
$ws.Range("R3").Value = 16.142942965765
$ws.Range("S3").Value = 0.00255675243683218
$ws.Range("T3").Value = 0.002311830395072382

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cx3cl1"
$ws.Range("C4").Value = "Cx3cr1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 16.6160005
$ws.Range("H4").Value = 33.232001
$ws.Range("I4").Value = 0.6672963354196896
$ws.Range("J4").Value = 0.6022128312718646
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.83935533333333
$ws.Range("N4").Value = 71.518066
$ws.Range("O4").Value = 0.5641043909871933
$ws.Range("P4").Value = 0.5651913495915012
$ws.Range("Q4").Value = 396.1147401383443
$ws.Range("R4").Value = 2376.688440830066
$ws.Range("S4").Value = 0.3764247928999099
$ws.Range("T4").Value = 0.3403654828478642

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cx3cl1"
$ws.Range("C5").Value = "Cx3cr1"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 16.6160005
$ws.Range("H5").Value = 33.232001
$ws.Range("I5").Value = 0.6672963354196896
$ws.Range("J5").Value = 0.6022128312718646
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.474532666666666
$ws.Range("N5").Value = 13.423598
$ws.Range("O5").Value = 0.1058796888417942
$ws.Range("P5").Value = 0.1060837057589585
$ws.Range("Q5").Value = 74.34883702659964
$ws.Range("R5").Value = 446.0930221595979
$ws.Range("S5").Value = 0.0706531283595063
$ws.Range("T5").Value = 0.06388496879691381

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Cx3cl1"
$ws.Range("C6").Value = "Cx3cr1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.6160005
$ws.Range("H6").Value = 33.232001
$ws.Range("I6").Value = 0.6672963354196896
$ws.Range("J6").Value = 0.6022128312718646
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.54090633333333
$ws.Range("N6").Value = 40.622719
$ws.Range("O6").Value = 0.3204149027427403
$ws.Range("P6").Value = 0.3210323021834275
$ws.Range("Q6").Value = 224.9957064051198
$ws.Range("R6").Value = 1349.974238430719
$ws.Range("S6").Value = 0.2138116904140868
$ws.Range("T6").Value = 0.1933297716276067

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cx3cl1"
$ws.Range("C7").Value = "Cx3cr1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.382185666666667
$ws.Range("H7").Value = 16.146557
$ws.Range("I7").Value = 0.216147849292316
$ws.Range("J7").Value = 0.2925994076090256
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.2438225
$ws.Range("N7").Value = 0.487645
$ws.Range("O7").Value = 0.005769507646004085
$ws.Range("P7").Value = 0.003853749843732457
$ws.Range("Q7").Value = 1.312297964710833
$ws.Range("R7").Value = 7.873787788265001
$ws.Range("S7").Value = 0.001247066669159356
$ws.Range("T7").Value = 0.001127604921349492

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Cx3cl1"
$ws.Range("C8").Value = "Cx3cr1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.382185666666667
$ws.Range("H8").Value = 16.146557
$ws.Range("I8").Value = 0.216147849292316
$ws.Range("J8").Value = 0.2925994076090256
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1619216666666667
$ws.Range("N8").Value = 0.485765
$ws.Range("O8").Value = 0.003831509782268077
$ws.Range("P8").Value = 0.003838892622380414
$ws.Range("Q8").Value = 0.8714924734561111
$ws.Range("R8").Value = 7.843432261105001
$ws.Range("S8").Value = 0.0008281725989797145
$ws.Range("T8").Value = 0.001123257707183168

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Cx3cl1"
$ws.Range("C9").Value = "Cx3cr1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.382185666666667
$ws.Range("H9").Value = 16.146557
$ws.Range("I9").Value = 0.216147849292316
$ws.Range("J9").Value = 0.2925994076090256
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 23.83935533333333
$ws.Range("N9").Value = 71.518066
$ws.Range("O9").Value = 0.5641043909871933
$ws.Range("P9").Value = 0.5651913495915012
$ws.Range("Q9").Value = 128.3078365776402
$ws.Range("R9").Value = 1154.770529198762
$ws.Range("S9").Value = 0.1219299508882335
$ws.Range("T9").Value = 0.165374654076219

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Cx3cl1"
$ws.Range("C10").Value = "Cx3cr1"
$ws.Range("D10").Value = "Neutrophils"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.382185666666667
$ws.Range("H10").Value = 16.146557
$ws.Range("I10").Value = 0.216147849292316
$ws.Range("J10").Value = 0.2925994076090256
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.474532666666666
$ws.Range("N10").Value = 13.423598
$ws.Range("O10").Value = 0.1058796888417942
$ws.Range("P10").Value = 0.1060837057589585
$ws.Range("Q10").Value = 24.08276558356511
$ws.Range("R10").Value = 216.744890252086
$ws.Range("S10").Value = 0.02288566702689345
$ws.Range("T10").Value = 0.03104002946204144

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Cx3cl1"
$ws.Range("C11").Value = "Cx3cr1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.382185666666667
$ws.Range("H11").Value = 16.146557
$ws.Range("I11").Value = 0.216147849292316
$ws.Range("J11").Value = 0.2925994076090256
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.54090633333333
$ws.Range("N11").Value = 40.622719
$ws.Range("O11").Value = 0.3204149027427403
$ws.Range("P11").Value = 0.3210323021834275
$ws.Range("Q11").Value = 72.87967198094256
$ws.Range("R11").Value = 655.917047828483
$ws.Range("S11").Value = 0.06925699210904991
$ws.Range("T11").Value = 0.09393386144223259

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Cx3cl1"
$ws.Range("C12").Value = "Cx3cr1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.902296
$ws.Range("H12").Value = 5.804592
$ws.Range("I12").Value = 0.1165558152879945
$ws.Range("J12").Value = 0.1051877611191097
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.2438225
$ws.Range("N12").Value = 0.487645
$ws.Range("O12").Value = 0.005769507646004085
$ws.Range("P12").Value = 0.003853749843732457
$ws.Range("Q12").Value = 0.7076450664599999
$ws.Range("R12").Value = 2.83058026584
$ws.Range("S12").Value = 0.0006724696674903238
$ws.Range("T12").Value = 0.0004053673179753362

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Cx3cl1"
$ws.Range("C13").Value = "Cx3cr1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.902296
$ws.Range("H13").Value = 5.804592
$ws.Range("I13").Value = 0.1165558152879945
$ws.Range("J13").Value = 0.1051877611191097
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1619216666666667
$ws.Range("N13").Value = 0.485765
$ws.Range("O13").Value = 0.003831509782268077
$ws.Range("P13").Value = 0.003838892622380414
$ws.Range("Q13").Value = 0.4699446054799999
$ws.Range("R13").Value = 2.81966763288
$ws.Range("S13").Value = 0.0004465847464561818
$ws.Range("T13").Value = 0.0004038045201248637

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Cx3cl1"
$ws.Range("C14").Value = "Cx3cr1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.902296
$ws.Range("H14").Value = 5.804592
$ws.Range("I14").Value = 0.1165558152879945
$ws.Range("J14").Value = 0.1051877611191097
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 23.83935533333333
$ws.Range("N14").Value = 71.518066
$ws.Range("O14").Value = 0.5641043909871933
$ws.Range("P14").Value = 0.5651913495915012
$ws.Range("Q14").Value = 69.188865626512
$ws.Range("R14").Value = 415.133193759072
$ws.Range("S14").Value = 0.06574964719904991
$ws.Range("T14").Value = 0.05945121266741806

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Cx3cl1"
$ws.Range("C15").Value = "Cx3cr1"
$ws.Range("D15").Value = "Neutrophils"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.902296
$ws.Range("H15").Value = 5.804592
$ws.Range("I15").Value = 0.1165558152879945
$ws.Range("J15").Value = 0.1051877611191097
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.474532666666666
$ws.Range("N15").Value = 13.423598
$ws.Range("O15").Value = 0.1058796888417942
$ws.Range("P15").Value = 0.1060837057589585
$ws.Range("Q15").Value = 12.986418260336
$ws.Range("R15").Value = 77.91850956201598
$ws.Range("S15").Value = 0.0123408934553945
$ws.Range("T15").Value = 0.01115870750000325

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Cx3cl1"
$ws.Range("C16").Value = "Cx3cr1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.902296
$ws.Range("H16").Value = 5.804592
$ws.Range("I16").Value = 0.1165558152879945
$ws.Range("J16").Value = 0.1051877611191097
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 13.54090633333333
$ws.Range("N16").Value = 40.622719
$ws.Range("O16").Value = 0.3204149027427403
$ws.Range("P16").Value = 0.3210323021834275
$ws.Range("Q16").Value = 39.29971828760799
$ws.Range("R16").Value = 235.798309725648
$ws.Range("S16").Value = 0.03734622021960355
$ws.Range("T16").Value = 0.03376866911358822

# Remove the now-obsolete row 17 (Neutrophils sending-cluster rows no longer present)
$ws.Rows(17).Delete()
